# "Git Setup brushup.docx" -- commit "added docker and ansible"
#
# Three content changes, replayed via the Word object model:
#   1. The `git checkout ` run right after the "####Creating a release
#      version " heading (in the "Creating a release version" snippet)
#      turns from color 24292E to C00000.
#   2. The text "HotFix Branch " (inside "####Creating HotFix Branch ")
#      gets wrapped in a `_GoBack` bookmark.
#   3. Because a document can only have one bookmark named `_GoBack`,
#      adding it in its new spot automatically pulls it off the empty
#      paragraph at the end of the document, leaving that paragraph
#      with no content.

$d = $word.ActiveDocument

# --- 1. Recolor "git checkout " under "####Creating a release version" ---
$rColor = $d.Content
$found = $rColor.Find.Execute("####Creating a release version", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not locate '####Creating a release version' heading text" }
$rColor.Collapse(0)   # wdCollapseEnd - continue searching right after the heading run

$found = $rColor.Find.Execute("git checkout ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not locate the 'git checkout ' run to recolor" }
$rColor.Font.Color = 192   # BGR(0000C0) == RGB "C00000"

# --- 2. Move the _GoBack bookmark onto "HotFix Branch " ------------------
$rBookmark = $d.Content
$found = $rBookmark.Find.Execute("HotFix Branch ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not locate 'HotFix Branch ' text for the bookmark" }

# Re-adding a bookmark with an existing name moves it; this both creates
# the new bookmark around "HotFix Branch " and removes it from the empty
# trailing paragraph where it used to live.
$d.Bookmarks.Add("_GoBack", $rBookmark)
